$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "72.849.92"
$ws.Range("E2").Value = "  +1.89%  "
$ws.Range("D3").Value = "3.988.99"
$ws.Range("E3").Value = "  +0.22%  "
$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").Value = "'592.41"
$ws.Range("E5").Value = "  +9.38%  "
$ws.Range("D6").Value = "'160.36"
$ws.Range("E6").Value = "  +7.71%  "
$ws.Range("D7").Value = "'0.683"
$ws.Range("E7").Value = "  -0.62%  "
$ws.Range("D8").Value = "'0.998"
$ws.Range("E8").Value = "  -0.14%  "
$ws.Range("E9").Value = "  +0.70%  "
$ws.Range("E10").Value = "  +0.60%  "
$ws.Range("D11").Value = "'53.70"
$ws.Range("E12").Value = "  -0.41%  "
$ws.Range("D13").Value = "'10.95"
$ws.Range("E13").Value = "  +1.94%  "
$ws.Range("D14").Value = "4.620.44"
$ws.Range("E14").Value = "  +0.17%  "
$ws.Range("D15").Value = "3.998.03"
$ws.Range("E15").Value = "  +0.16%  "
$ws.Range("D16").Value = "'1.27"
$ws.Range("E16").Value = "  +8.39%  "
$ws.Range("D18").Value = "'20.36"
$ws.Range("E18").Value = "  -1.15%  "
$ws.Range("E19").Value = "  +0.17%  "
$ws.Range("D20").Value = "72.549.43"
$ws.Range("E20").Value = "  +1.66%  "
$ws.Range("D21").Value = "'433.06"
$ws.Range("E21").Value = "  +0.99%  "
$ws.Range("D22").Value = "'4.78"
$ws.Range("E22").Value = "  +13.40%  "
$ws.Range("D23").Value = "'96.19"
$ws.Range("E23").Value = "  -1.37%  "
$ws.Range("E24").Value = "  -4.66%  "
$ws.Range("D25").Value = "'14.17"
$ws.Range("E25").Value = "  -2.77%  "
$ws.Range("D26").Value = "'4.41"
$ws.Range("E26").Value = "  +17.72%  "
$ws.Range("E27").Value = "  -2.04%  "
$ws.Range("E28").Value = "  +0.96%  "
$ws.Range("E29").Value = "  -2.66%  "
$ws.Range("D30").Value = "'36.35"
$ws.Range("E30").Value = "  -0.88%  "
$ws.Range("E31").Value = "  +2.70%  "
$ws.Range("D32").Value = "'13.77"
$ws.Range("E32").Value = "  +2.49%  "
$ws.Range("E33").Value = "  -0.42%  "
$ws.Range("D34").Value = "'48.82"
$ws.Range("E34").Value = "  -4.86%  "
$ws.Range("D35").Value = "'670.23"
$ws.Range("E35").Value = "  -3.49%  "
$ws.Range("D36").Value = "'70.41"
$ws.Range("E36").Value = "  +7.43%  "
$ws.Range("D37").Value = "'0.438"
$ws.Range("E37").Value = "  -0.11%  "
$ws.Range("D38").Value = "0.0₃0876"
$ws.Range("E38").Value = "  +6.53%  "
$ws.Range("D39").Value = "'0.146"
$ws.Range("E39").Value = "  -3.20%  "
$ws.Range("B40").Value = "Dai"
$ws.Range("C40").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D40").Value = "'0.999"
$ws.Range("E40").Value = "  -0.46%  "
$ws.Range("B41").Value = "ThetaToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D41").Value = "'3.34"
$ws.Range("E41").Value = "  -2.31%  "
$ws.Range("D42").Value = "'3.34"
$ws.Range("E42").Value = "  +1.58%  "
$ws.Range("E43").Value = "  +0.30%  "
$ws.Range("E44").Value = "  +0.83%  "
$ws.Range("D45").Value = "'10.68"
$ws.Range("E45").Value = "  +9.46%  "
$ws.Range("E46").Value = "  +0.36%  "
$ws.Range("D47").Value = "'3.46"
$ws.Range("E47").Value = "  +3.19%  "
$ws.Range("E48").Value = "  -4.18%  "
$ws.Range("D49").Value = "2.852.50"
$ws.Range("E49").Value = "  +1.41%  "
$ws.Range("D50").Value = "'3.02"
$ws.Range("E50").Value = "  +0.54%  "
$ws.Range("E51").Value = "  +4.03%  "
